$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching style of existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column (H2:H24)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 1
    24 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
